$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1831944722881327
$ws.Range("D2").Value = 0.04429341463183867
$ws.Range("E2").Value = 0.2603791776828435

$ws.Range("C3").Value = 0.07354921245638879
$ws.Range("D3").Value = 0.1451621806139535
$ws.Range("E3").Value = 0.3378236949304282

$ws.Range("C4").Value = 0.5140503535480559
$ws.Range("D4").Value = 0.7738762166124983
$ws.Range("E4").Value = 0.1778277411902804

$ws.Range("C5").Value = 0.04812471609263469
$ws.Range("D5").Value = 0.4372860710202552
$ws.Range("E5").Value = 0.0266144070951404

$ws.Range("C6").Value = 0.1622413097994405
$ws.Range("D6").Value = 0.8164297880401218
$ws.Range("E6").Value = 0.2471189745357338

$ws.Range("C7").Value = 0.1154468374321959
$ws.Range("D7").Value = 0.3985122619810699
$ws.Range("E7").Value = 0.1252276747800244
